# Auto-generated PowerShell COM-interop script
# Applies the sweep_measurements.xlsx update (commit: "20250923 project is
# working well, commit to save a know good"):
#   - Measurements sheet: refresh the sweep data for rows 2-3 and append
#     5 newly captured sweep points as rows 4-8 (dimension grows to A1:AX8)
#   - Row heights normalized from the old 60pt auto-wrap height to 15pt
#   - The per-row free-text Comment (column AX) is left blank going forward
#   - Statistics sheet: Max/Min/Mean rollups recomputed for the new dataset

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Measurements")
$stats = $wb.Worksheets.Item("Statistics")

# Row 2
$ws.Rows.Item(2).RowHeight = 15
$ws.Cells.Item(2, 1).Value = 1.220230579376221
$ws.Cells.Item(2, 2).Value = 3.601427555084229
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = 6
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.555
$ws.Cells.Item(2, 7).Value = -10.99513
$ws.Cells.Item(2, 8).Value = 5.929387
$ws.Cells.Item(2, 9).Value = 5.905197144
$ws.Cells.Item(2, 10).Value = -43.7933197
$ws.Cells.Item(2, 11).Value = 0.1043274402618408
$ws.Cells.Item(2, 12).Value = 5.93272781372
$ws.Cells.Item(2, 13).Value = -51.6141395569
$ws.Cells.Item(2, 14).Value = -53.0046730042
$ws.Cells.Item(2, 15).Value = 0.1998167037963867
$ws.Cells.Item(2, 16).Value = 5.971302032
$ws.Cells.Item(2, 17).Value = -50.41112518
$ws.Cells.Item(2, 18).Value = 0.0866706371307373
$ws.Cells.Item(2, 19).Value = 6.04659557343
$ws.Cells.Item(2, 20).Value = -49.6169271469
$ws.Cells.Item(2, 21).Value = -52.4700574875
$ws.Cells.Item(2, 22).Value = 0.1782207489013672
$ws.Cells.Item(2, 23).Value = 6.827713489532471
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0.409
$ws.Cells.Item(2, 27).Value = 5.982486725
$ws.Cells.Item(2, 28).Value = -49.02235031
$ws.Cells.Item(2, 29).Value = 0.15743088722229
$ws.Cells.Item(2, 30).Value = 6.05910825729
$ws.Cells.Item(2, 31).Value = -48.4489855766
$ws.Cells.Item(2, 32).Value = -50.3370471001
$ws.Cells.Item(2, 33).Value = 0.1763713359832764
$ws.Cells.Item(2, 34).Value = 2.707039833068848
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 0
$ws.Cells.Item(2, 37).Value = 0.4
$ws.Cells.Item(2, 38).Value = 5.914958954
$ws.Cells.Item(2, 39).Value = -49.59848022
$ws.Cells.Item(2, 40).Value = 0.1543612480163574
$ws.Cells.Item(2, 41).Value = 5.94747495651
$ws.Cells.Item(2, 42).Value = -59.9518580437
$ws.Cells.Item(2, 43).Value = -60.6399531364
$ws.Cells.Item(2, 44).Value = 0.1940133571624756
$ws.Cells.Item(2, 45).Value = 9.742092847824097
$ws.Cells.Item(2, 46).Value = 0
$ws.Cells.Item(2, 47).Value = 0
$ws.Cells.Item(2, 48).Value = 0.425
$ws.Cells.Item(2, 49).Value = 27.307
$ws.Range("AX2").Value = ""
$ws.Range("AX2").WrapText = $true

# Row 3
$ws.Rows.Item(3).RowHeight = 15
$ws.Cells.Item(3, 1).Value = 1.220230579376221
$ws.Cells.Item(3, 2).Value = 3.601427555084229
$ws.Cells.Item(3, 3).Value = 4.1
$ws.Cells.Item(3, 4).Value = 6
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.403
$ws.Cells.Item(3, 7).Value = -10.71387
$ws.Cells.Item(3, 8).Value = 6.003528
$ws.Cells.Item(3, 9).Value = 5.952327728
$ws.Cells.Item(3, 10).Value = -43.88660812
$ws.Cells.Item(3, 11).Value = 0.1090087890625
$ws.Cells.Item(3, 12).Value = 6.02062749863
$ws.Cells.Item(3, 13).Value = -51.8926329613
$ws.Cells.Item(3, 14).Value = -53.2699217796
$ws.Cells.Item(3, 15).Value = 0.192603588104248
$ws.Cells.Item(3, 16).Value = 5.913562775
$ws.Cells.Item(3, 17).Value = -50.61495209
$ws.Cells.Item(3, 18).Value = 0.09663581848144531
$ws.Cells.Item(3, 19).Value = 6.02402734756
$ws.Cells.Item(3, 20).Value = -50.4292626381
$ws.Cells.Item(3, 21).Value = -52.0804390907
$ws.Cells.Item(3, 22).Value = 0.1820697784423828
$ws.Cells.Item(3, 23).Value = 6.274100780487061
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(3, 25).Value = 0
$ws.Cells.Item(3, 26).Value = 0.181
$ws.Cells.Item(3, 27).Value = 5.911529541
$ws.Cells.Item(3, 28).Value = -48.9655838
$ws.Cells.Item(3, 29).Value = 0.1908462047576904
$ws.Cells.Item(3, 30).Value = 6.0308675766
$ws.Cells.Item(3, 31).Value = -48.5944280624
$ws.Cells.Item(3, 32).Value = -50.3751173019
$ws.Cells.Item(3, 33).Value = 0.1806857585906982
$ws.Cells.Item(3, 34).Value = 2.836046457290649
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = 0
$ws.Cells.Item(3, 37).Value = 0.43
$ws.Cells.Item(3, 38).Value = 5.951593399
$ws.Cells.Item(3, 39).Value = -49.17251587
$ws.Cells.Item(3, 40).Value = 0.1651663780212402
$ws.Cells.Item(3, 41).Value = 6.03072452545
$ws.Cells.Item(3, 42).Value = -59.9182367325
$ws.Cells.Item(3, 43).Value = -60.3616609573
$ws.Cells.Item(3, 44).Value = 0.1942462921142578
$ws.Cells.Item(3, 45).Value = 9.311590671539307
$ws.Cells.Item(3, 46).Value = 0
$ws.Cells.Item(3, 47).Value = 0
$ws.Cells.Item(3, 48).Value = 0.446
$ws.Cells.Item(3, 49).Value = 26.24
$ws.Range("AX3").Value = ""
$ws.Range("AX3").WrapText = $true

# Row 4
$ws.Rows.Item(4).RowHeight = 15
$ws.Cells.Item(4, 1).Value = 1.220230579376221
$ws.Cells.Item(4, 2).Value = 3.601427555084229
$ws.Cells.Item(4, 3).Value = 4.2
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1.2
$ws.Cells.Item(4, 7).Value = -10.82424
$ws.Cells.Item(4, 8).Value = 6.005121
$ws.Cells.Item(4, 9).Value = 5.904348373
$ws.Cells.Item(4, 10).Value = -44.29456329
$ws.Cells.Item(4, 11).Value = 0.07992815971374512
$ws.Cells.Item(4, 12).Value = 5.99493026733
$ws.Cells.Item(4, 13).Value = -52.0992126465
$ws.Cells.Item(4, 14).Value = -53.5028305054
$ws.Cells.Item(4, 15).Value = 0.186129093170166
$ws.Cells.Item(4, 16).Value = 5.870141983
$ws.Cells.Item(4, 17).Value = -50.57678986
$ws.Cells.Item(4, 18).Value = 0.07984685897827148
$ws.Cells.Item(4, 19).Value = 6.00778341293
$ws.Cells.Item(4, 20).Value = -49.603662014
$ws.Cells.Item(4, 21).Value = -52.3552794456
$ws.Cells.Item(4, 22).Value = 0.1545419692993164
$ws.Cells.Item(4, 23).Value = 7.262543916702271
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 0.188
$ws.Cells.Item(4, 27).Value = 5.975738525
$ws.Cells.Item(4, 28).Value = -49.47615814
$ws.Cells.Item(4, 29).Value = 0.1377959251403809
$ws.Cells.Item(4, 30).Value = 6.11211681366
$ws.Cells.Item(4, 31).Value = -48.6733617783
$ws.Cells.Item(4, 32).Value = -50.8653478622
$ws.Cells.Item(4, 33).Value = 0.1530272960662842
$ws.Cells.Item(4, 34).Value = 2.78815746307373
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 0
$ws.Cells.Item(4, 37).Value = 0.411
$ws.Cells.Item(4, 38).Value = 5.913640976
$ws.Cells.Item(4, 39).Value = -50.08543396
$ws.Cells.Item(4, 40).Value = 0.1237905025482178
$ws.Cells.Item(4, 41).Value = 5.99904155731
$ws.Cells.Item(4, 42).Value = -60.6964216232
$ws.Cells.Item(4, 43).Value = -59.9297742844
$ws.Cells.Item(4, 44).Value = 0.1895630359649658
$ws.Cells.Item(4, 45).Value = 9.392547369003296
$ws.Cells.Item(4, 46).Value = 0
$ws.Cells.Item(4, 47).Value = 0
$ws.Cells.Item(4, 48).Value = 0.403
$ws.Cells.Item(4, 49).Value = 27.075
$ws.Range("AX4").Value = ""
$ws.Range("AX4").WrapText = $true

# Row 5
$ws.Rows.Item(5).RowHeight = 15
$ws.Cells.Item(5, 1).Value = 1.220230579376221
$ws.Cells.Item(5, 2).Value = 3.601427555084229
$ws.Cells.Item(5, 3).Value = 4.3
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.606
$ws.Cells.Item(5, 7).Value = -10.96826
$ws.Cells.Item(5, 8).Value = 6.032585
$ws.Cells.Item(5, 9).Value = 5.959341049
$ws.Cells.Item(5, 10).Value = -44.20947647
$ws.Cells.Item(5, 11).Value = 0.09226107597351074
$ws.Cells.Item(5, 12).Value = 6.04576396942
$ws.Cells.Item(5, 13).Value = -51.904009819
$ws.Cells.Item(5, 14).Value = -53.4327001572
$ws.Cells.Item(5, 15).Value = 0.1979174613952637
$ws.Cells.Item(5, 16).Value = 5.921918869
$ws.Cells.Item(5, 17).Value = -50.83035278
$ws.Cells.Item(5, 18).Value = 0.04189252853393555
$ws.Cells.Item(5, 19).Value = 6.05037403107
$ws.Cells.Item(5, 20).Value = -49.6174402237
$ws.Cells.Item(5, 21).Value = -52.6674547195
$ws.Cells.Item(5, 22).Value = 0.1787078380584717
$ws.Cells.Item(5, 23).Value = 7.230472564697266
$ws.Cells.Item(5, 24).Value = 0
$ws.Cells.Item(5, 25).Value = 0
$ws.Cells.Item(5, 26).Value = 0.176
$ws.Cells.Item(5, 27).Value = 5.922689438
$ws.Cells.Item(5, 28).Value = -49.75502396
$ws.Cells.Item(5, 29).Value = 0.1813395023345947
$ws.Cells.Item(5, 30).Value = 6.05016994476
$ws.Cells.Item(5, 31).Value = -48.7619047165
$ws.Cells.Item(5, 32).Value = -50.5073003769
$ws.Cells.Item(5, 33).Value = 0.1777670383453369
$ws.Cells.Item(5, 34).Value = 2.725986242294312
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 0
$ws.Cells.Item(5, 37).Value = 0.384
$ws.Cells.Item(5, 38).Value = 5.953754425
$ws.Cells.Item(5, 39).Value = -50.27392578
$ws.Cells.Item(5, 40).Value = 0.1299064159393311
$ws.Cells.Item(5, 41).Value = 6.04695129395
$ws.Cells.Item(5, 42).Value = -60.7487335205
$ws.Cells.Item(5, 43).Value = -61.0186157227
$ws.Cells.Item(5, 44).Value = 0.1886463165283203
$ws.Cells.Item(5, 45).Value = 9.319432497024536
$ws.Cells.Item(5, 46).Value = 0
$ws.Cells.Item(5, 47).Value = 0
$ws.Cells.Item(5, 48).Value = 0.423
$ws.Cells.Item(5, 49).Value = 27.224
$ws.Range("AX5").Value = ""
$ws.Range("AX5").WrapText = $true

# Row 6
$ws.Rows.Item(6).RowHeight = 15
$ws.Cells.Item(6, 1).Value = 1.220230579376221
$ws.Cells.Item(6, 2).Value = 3.601427555084229
$ws.Cells.Item(6, 3).Value = 4.4
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.594
$ws.Cells.Item(6, 7).Value = -10.93479
$ws.Cells.Item(6, 8).Value = 5.937651
$ws.Cells.Item(6, 9).Value = 5.923822403
$ws.Cells.Item(6, 10).Value = -43.79166031
$ws.Cells.Item(6, 11).Value = 0.1474802494049072
$ws.Cells.Item(6, 12).Value = 5.96760511398
$ws.Cells.Item(6, 13).Value = -51.4035868645
$ws.Cells.Item(6, 14).Value = -52.8622965813
$ws.Cells.Item(6, 15).Value = 0.1923964023590088
$ws.Cells.Item(6, 16).Value = 5.884895325
$ws.Cells.Item(6, 17).Value = -50.72265244
$ws.Cells.Item(6, 18).Value = 0.08623671531677246
$ws.Cells.Item(6, 19).Value = 5.97079944611
$ws.Cells.Item(6, 20).Value = -49.2488260269
$ws.Cells.Item(6, 21).Value = -52.2822999954
$ws.Cells.Item(6, 22).Value = 0.1858205795288086
$ws.Cells.Item(6, 23).Value = 7.333505630493164
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(6, 26).Value = 0.195
$ws.Cells.Item(6, 27).Value = 5.985588074
$ws.Cells.Item(6, 28).Value = -49.01224136
$ws.Cells.Item(6, 29).Value = 0.1716275215148926
$ws.Cells.Item(6, 30).Value = 6.06830358505
$ws.Cells.Item(6, 31).Value = -48.481107235
$ws.Cells.Item(6, 32).Value = -50.1791052818
$ws.Cells.Item(6, 33).Value = 0.1749708652496338
$ws.Cells.Item(6, 34).Value = 2.733812808990479
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = 0
$ws.Cells.Item(6, 37).Value = 0.406
$ws.Cells.Item(6, 38).Value = 5.921880722
$ws.Cells.Item(6, 39).Value = -49.69962311
$ws.Cells.Item(6, 40).Value = 0.1499216556549072
$ws.Cells.Item(6, 41).Value = 5.96842288971
$ws.Cells.Item(6, 42).Value = -59.8529634476
$ws.Cells.Item(6, 43).Value = -60.2399606705
$ws.Cells.Item(6, 44).Value = 0.1987082958221436
$ws.Cells.Item(6, 45).Value = 9.361552000045776
$ws.Cells.Item(6, 46).Value = 0
$ws.Cells.Item(6, 47).Value = 0
$ws.Cells.Item(6, 48).Value = 0.446
$ws.Cells.Item(6, 49).Value = 27.281
$ws.Range("AX6").Value = ""
$ws.Range("AX6").WrapText = $true

# Row 7
$ws.Rows.Item(7).RowHeight = 15
$ws.Cells.Item(7, 1).Value = 1.220230579376221
$ws.Cells.Item(7, 2).Value = 3.601427555084229
$ws.Cells.Item(7, 3).Value = 4.5
$ws.Cells.Item(7, 4).Value = 6
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1.192
$ws.Cells.Item(7, 7).Value = -10.7707
$ws.Cells.Item(7, 8).Value = 6.005444
$ws.Cells.Item(7, 9).Value = 5.838497162
$ws.Cells.Item(7, 10).Value = -43.88553238
$ws.Cells.Item(7, 11).Value = 0.1202778816223145
$ws.Cells.Item(7, 12).Value = 5.91196155548
$ws.Cells.Item(7, 13).Value = -51.5794992447
$ws.Cells.Item(7, 14).Value = -52.9670495987
$ws.Cells.Item(7, 15).Value = 0.1932497024536133
$ws.Cells.Item(7, 16).Value = 5.791469574
$ws.Cells.Item(7, 17).Value = -50.85944366
$ws.Cells.Item(7, 18).Value = 0.1074469089508057
$ws.Cells.Item(7, 19).Value = 5.89789676666
$ws.Cells.Item(7, 20).Value = -50.125210762
$ws.Cells.Item(7, 21).Value = -52.4777765274
$ws.Cells.Item(7, 22).Value = 0.1501402854919434
$ws.Cells.Item(7, 23).Value = 6.274209260940552
$ws.Cells.Item(7, 24).Value = 0
$ws.Cells.Item(7, 25).Value = 0
$ws.Cells.Item(7, 26).Value = 0.186
$ws.Cells.Item(7, 27).Value = 5.887290955
$ws.Cells.Item(7, 28).Value = -49.04273224
$ws.Cells.Item(7, 29).Value = 0.1730093955993652
$ws.Cells.Item(7, 30).Value = 6.00402545929
$ws.Cells.Item(7, 31).Value = -48.8432817459
$ws.Cells.Item(7, 32).Value = -50.3793001175
$ws.Cells.Item(7, 33).Value = 0.1541531085968018
$ws.Cells.Item(7, 34).Value = 2.798098087310791
$ws.Cells.Item(7, 35).Value = 0
$ws.Cells.Item(7, 36).Value = 0
$ws.Cells.Item(7, 37).Value = 0.436
$ws.Cells.Item(7, 38).Value = 5.92755127
$ws.Cells.Item(7, 39).Value = -49.1100502
$ws.Cells.Item(7, 40).Value = 0.1650357246398926
$ws.Cells.Item(7, 41).Value = 6.00565910339
$ws.Cells.Item(7, 42).Value = -59.1627941132
$ws.Cells.Item(7, 43).Value = -59.5720386505
$ws.Cells.Item(7, 44).Value = 0.1917657852172852
$ws.Cells.Item(7, 45).Value = 9.409790992736816
$ws.Cells.Item(7, 46).Value = 0
$ws.Cells.Item(7, 47).Value = 0
$ws.Cells.Item(7, 48).Value = 0.432
$ws.Cells.Item(7, 49).Value = 26.005
$ws.Range("AX7").Value = ""
$ws.Range("AX7").WrapText = $true

# Row 8
$ws.Rows.Item(8).RowHeight = 15
$ws.Cells.Item(8, 1).Value = 1.220230579376221
$ws.Cells.Item(8, 2).Value = 3.601427555084229
$ws.Cells.Item(8, 3).Value = 4.6
$ws.Cells.Item(8, 4).Value = 6
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1.191
$ws.Cells.Item(8, 7).Value = -10.74574
$ws.Cells.Item(8, 8).Value = 6.011664
$ws.Cells.Item(8, 9).Value = 5.96836853
$ws.Cells.Item(8, 10).Value = -43.80942535
$ws.Cells.Item(8, 11).Value = 0.0797128677368164
$ws.Cells.Item(8, 12).Value = 6.04242086411
$ws.Cells.Item(8, 13).Value = -51.6768660545
$ws.Cells.Item(8, 14).Value = -53.306283474
$ws.Cells.Item(8, 15).Value = 0.1866886615753174
$ws.Cells.Item(8, 16).Value = 5.921043396
$ws.Cells.Item(8, 17).Value = -50.95992661
$ws.Cells.Item(8, 18).Value = 0.08726811408996582
$ws.Cells.Item(8, 19).Value = 6.03580713272
$ws.Cells.Item(8, 20).Value = -49.1943421364
$ws.Cells.Item(8, 21).Value = -52.8759865761
$ws.Cells.Item(8, 22).Value = 0.1773674488067627
$ws.Cells.Item(8, 23).Value = 6.246154546737671
$ws.Cells.Item(8, 24).Value = 0
$ws.Cells.Item(8, 25).Value = 0
$ws.Cells.Item(8, 26).Value = 0.2
$ws.Cells.Item(8, 27).Value = 5.918746948
$ws.Cells.Item(8, 28).Value = -48.86702728
$ws.Cells.Item(8, 29).Value = 0.2074992656707764
$ws.Cells.Item(8, 30).Value = 6.03691387177
$ws.Cells.Item(8, 31).Value = -47.9476156235
$ws.Cells.Item(8, 32).Value = -50.5188627243
$ws.Cells.Item(8, 33).Value = 0.1777260303497314
$ws.Cells.Item(8, 34).Value = 2.776440382003784
$ws.Cells.Item(8, 35).Value = 0
$ws.Cells.Item(8, 36).Value = 0
$ws.Cells.Item(8, 37).Value = 0.413
$ws.Cells.Item(8, 38).Value = 5.953157425
$ws.Cells.Item(8, 39).Value = -49.09585953
$ws.Cells.Item(8, 40).Value = 0.1705596446990967
$ws.Cells.Item(8, 41).Value = 6.03079891205
$ws.Cells.Item(8, 42).Value = -59.3464918137
$ws.Cells.Item(8, 43).Value = -59.8654317856
$ws.Cells.Item(8, 44).Value = 0.1938447952270508
$ws.Cells.Item(8, 45).Value = 9.272533416748047
$ws.Cells.Item(8, 46).Value = 0
$ws.Cells.Item(8, 47).Value = 0
$ws.Cells.Item(8, 48).Value = 0.399
$ws.Cells.Item(8, 49).Value = 26.135
$ws.Range("AX8").Value = ""
$ws.Range("AX8").WrapText = $true

# Statistics sheet: recompute Max/Min/Mean rows for the refreshed dataset
$stats.Cells.Item(2, 2).Value = 7
$stats.Cells.Item(3, 2).Value = 1.220230579376221
$stats.Cells.Item(4, 2).Value = 1.220230579376221
$stats.Cells.Item(5, 2).Value = 1.220230579376221
$stats.Cells.Item(6, 2).Value = 3.601427555084229
$stats.Cells.Item(7, 2).Value = 3.601427555084229
$stats.Cells.Item(8, 2).Value = 3.601427555084229
$stats.Cells.Item(9, 2).Value = 4.6
$stats.Cells.Item(10, 2).Value = 4
$stats.Cells.Item(11, 2).Value = 4.3
$stats.Cells.Item(12, 2).Value = 6
$stats.Cells.Item(13, 2).Value = 6
$stats.Cells.Item(14, 2).Value = 6
$stats.Cells.Item(15, 2).Value = 2
$stats.Cells.Item(16, 2).Value = 1
$stats.Cells.Item(17, 2).Value = 1.4285714285714286
$stats.Cells.Item(18, 2).Value = 1.2
$stats.Cells.Item(19, 2).Value = 0.403
$stats.Cells.Item(20, 2).Value = 0.8201428571428571
$stats.Cells.Item(21, 2).Value = -10.71387
$stats.Cells.Item(22, 2).Value = -10.99513
$stats.Cells.Item(23, 2).Value = -10.85039
$stats.Cells.Item(24, 2).Value = 6.032585
$stats.Cells.Item(25, 2).Value = 5.929387
$stats.Cells.Item(26, 2).Value = 5.989339999999999
$stats.Cells.Item(27, 2).Value = 5.96836853
$stats.Cells.Item(28, 2).Value = 5.838497162
$stats.Cells.Item(29, 2).Value = 5.921700341285714
$stats.Cells.Item(30, 2).Value = -43.79166031
$stats.Cells.Item(31, 2).Value = -44.29456329
$stats.Cells.Item(32, 2).Value = -43.95294080285714
$stats.Cells.Item(33, 2).Value = 0.1474802494049072
$stats.Cells.Item(34, 2).Value = 0.0797128677368164
$stats.Cells.Item(35, 2).Value = 0.1047137805393764
$stats.Cells.Item(36, 2).Value = 6.04576396942
$stats.Cells.Item(37, 2).Value = 5.91196155548
$stats.Cells.Item(38, 2).Value = 5.988005297524286
$stats.Cells.Item(39, 2).Value = -51.4035868645
$stats.Cells.Item(40, 2).Value = -52.0992126465
$stats.Cells.Item(41, 2).Value = -51.738563878200004
$stats.Cells.Item(42, 2).Value = -52.8622965813
$stats.Cells.Item(43, 2).Value = -53.5028305054
$stats.Cells.Item(44, 2).Value = -53.192250728628565
$stats.Cells.Item(45, 2).Value = 0.1998167037963867
$stats.Cells.Item(46, 2).Value = 0.186129093170166
$stats.Cells.Item(47, 2).Value = 0.19268594469342912
$stats.Cells.Item(48, 2).Value = 5.971302032
$stats.Cells.Item(49, 2).Value = 5.791469574
$stats.Cells.Item(50, 2).Value = 5.896333422
$stats.Cells.Item(51, 2).Value = -50.41112518
$stats.Cells.Item(52, 2).Value = -50.95992661
$stats.Cells.Item(53, 2).Value = -50.710748945714286
$stats.Cells.Item(54, 2).Value = 0.1074469089508057
$stats.Cells.Item(55, 2).Value = 0.04189252853393555
$stats.Cells.Item(56, 2).Value = 0.0837139402117048
$stats.Cells.Item(57, 2).Value = 6.05037403107
$stats.Cells.Item(58, 2).Value = 5.89789676666
$stats.Cells.Item(59, 2).Value = 6.004754815782857
$stats.Cells.Item(60, 2).Value = -49.1943421364
$stats.Cells.Item(61, 2).Value = -50.4292626381
$stats.Cells.Item(62, 2).Value = -49.69081013542858
$stats.Cells.Item(63, 2).Value = -52.0804390907
$stats.Cells.Item(64, 2).Value = -52.8759865761
$stats.Cells.Item(65, 2).Value = -52.458470548885714
$stats.Cells.Item(66, 2).Value = 0.1858205795288086
$stats.Cells.Item(67, 2).Value = 0.1501402854919434
$stats.Cells.Item(68, 2).Value = 0.17240980693272182
$stats.Cells.Item(69, 2).Value = 7.333505630493164
$stats.Cells.Item(70, 2).Value = 6.246154546737671
$stats.Cells.Item(71, 2).Value = 6.778385741370065
$stats.Cells.Item(72, 2).Value = 0
$stats.Cells.Item(73, 2).Value = 0
$stats.Cells.Item(74, 2).Value = 0
$stats.Cells.Item(75, 2).Value = 0
$stats.Cells.Item(76, 2).Value = 0
$stats.Cells.Item(77, 2).Value = 0
$stats.Cells.Item(78, 2).Value = 0.409
$stats.Cells.Item(79, 2).Value = 0.176
$stats.Cells.Item(80, 2).Value = 0.21928571428571428
$stats.Cells.Item(81, 2).Value = 5.985588074
$stats.Cells.Item(82, 2).Value = 5.887290955
$stats.Cells.Item(83, 2).Value = 5.940581458
$stats.Cells.Item(84, 2).Value = -48.86702728
$stats.Cells.Item(85, 2).Value = -49.75502396
$stats.Cells.Item(86, 2).Value = -49.16301672714286
$stats.Cells.Item(87, 2).Value = 0.2074992656707764
$stats.Cells.Item(88, 2).Value = 0.1377959251403809
$stats.Cells.Item(89, 2).Value = 0.17422124317714147
$stats.Cells.Item(90, 2).Value = 6.11211681366
$stats.Cells.Item(91, 2).Value = 6.00402545929
$stats.Cells.Item(92, 2).Value = 6.05164364406
$stats.Cells.Item(93, 2).Value = -47.9476156235
$stats.Cells.Item(94, 2).Value = -48.8432817459
$stats.Cells.Item(95, 2).Value = -48.53581210545714
$stats.Cells.Item(96, 2).Value = -50.1791052818
$stats.Cells.Item(97, 2).Value = -50.8653478622
$stats.Cells.Item(98, 2).Value = -50.45172582352858
$stats.Cells.Item(99, 2).Value = 0.1806857585906982
$stats.Cells.Item(100, 2).Value = 0.1530272960662842
$stats.Cells.Item(101, 2).Value = 0.1706716333116804
$stats.Cells.Item(102, 2).Value = 2.836046457290649
$stats.Cells.Item(103, 2).Value = 2.707039833068848
$stats.Cells.Item(104, 2).Value = 2.7665116105760847
$stats.Cells.Item(105, 2).Value = 0
$stats.Cells.Item(106, 2).Value = 0
$stats.Cells.Item(107, 2).Value = 0
$stats.Cells.Item(108, 2).Value = 0
$stats.Cells.Item(109, 2).Value = 0
$stats.Cells.Item(110, 2).Value = 0
$stats.Cells.Item(111, 2).Value = 0.436
$stats.Cells.Item(112, 2).Value = 0.384
$stats.Cells.Item(113, 2).Value = 0.4114285714285714
$stats.Cells.Item(114, 2).Value = 5.953754425
$stats.Cells.Item(115, 2).Value = 5.913640976
$stats.Cells.Item(116, 2).Value = 5.933791024428571
$stats.Cells.Item(117, 2).Value = -49.09585953
$stats.Cells.Item(118, 2).Value = -50.27392578
$stats.Cells.Item(119, 2).Value = -49.576555524285716
$stats.Cells.Item(120, 2).Value = 0.1705596446990967
$stats.Cells.Item(121, 2).Value = 0.1237905025482178
$stats.Cells.Item(122, 2).Value = 0.15124879564557756
$stats.Cells.Item(123, 2).Value = 6.04695129395
$stats.Cells.Item(124, 2).Value = 5.94747495651
$stats.Cells.Item(125, 2).Value = 6.004153319767143
$stats.Cells.Item(126, 2).Value = -59.1627941132
$stats.Cells.Item(127, 2).Value = -60.7487335205
$stats.Cells.Item(128, 2).Value = -59.95392847062857
$stats.Cells.Item(129, 2).Value = -59.5720386505
$stats.Cells.Item(130, 2).Value = -61.0186157227
$stats.Cells.Item(131, 2).Value = -60.23249074391429
$stats.Cells.Item(132, 2).Value = 0.1987082958221436
$stats.Cells.Item(133, 2).Value = 0.1886463165283203
$stats.Cells.Item(134, 2).Value = 0.19296969686235704
$stats.Cells.Item(135, 2).Value = 9.742092847824097
$stats.Cells.Item(136, 2).Value = 9.272533416748047
$stats.Cells.Item(137, 2).Value = 9.401362827845983
$stats.Cells.Item(138, 2).Value = 0
$stats.Cells.Item(139, 2).Value = 0
$stats.Cells.Item(140, 2).Value = 0
$stats.Cells.Item(141, 2).Value = 0
$stats.Cells.Item(142, 2).Value = 0
$stats.Cells.Item(143, 2).Value = 0
$stats.Cells.Item(144, 2).Value = 0.446
$stats.Cells.Item(145, 2).Value = 0.399
$stats.Cells.Item(146, 2).Value = 0.4248571428571429
$stats.Cells.Item(147, 2).Value = 27.307
$stats.Cells.Item(148, 2).Value = 26.005
$stats.Cells.Item(149, 2).Value = 26.75242857142857

Write-Host "edit applied"
